# Auto-generated edit script: apply numeric corrections to Leve profit sheets
# per the commit diff (scheduled runner update pulling fresh market data).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 10.333333
$ws.Range("I8").Value = 10.333333
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 30.999999
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 108.000001
$ws.Range("N8").ClearContents()
$ws.Range("H40").Value = 8295.826999999999
$ws.Range("I40").Value = 8014.5625
$ws.Range("J40").Value = 8642
$ws.Range("K40").Value = 8014.5625
$ws.Range("L40").Value = 8642
$ws.Range("M40").Value = -7839.5625
$ws.Range("N40").Value = -8992
$ws.Range("H86").Value = 4546.077
$ws.Range("I86").Value = 3387.375
$ws.Range("K86").Value = 3387.375
$ws.Range("M86").Value = -2264.375
$ws.Range("H89").Value = 4546.077
$ws.Range("I89").Value = 3387.375
$ws.Range("K89").Value = 16936.875
$ws.Range("M89").Value = -11320.875
$ws.Range("H121").Value = 1777.9166
$ws.Range("J121").Value = 1777.9166
$ws.Range("L121").Value = 5333.7498
$ws.Range("N121").Value = -8827.7498
$ws.Range("H130").Value = 20000
$ws.Range("J130").Value = 20000
$ws.Range("L130").Value = 20000
$ws.Range("N130").Value = -30040
$ws.Range("H132").Value = 18575.719
$ws.Range("I132").Value = 1594
$ws.Range("K132").Value = 4782
$ws.Range("M132").Value = -2252
$ws.Range("H135").Value = 29421820
$ws.Range("I135").Value = 50006344
$ws.Range("K135").Value = 450057096
$ws.Range("M135").Value = -450054561
$ws.Range("H137").Value = 18207402
$ws.Range("I137").Value = 33375570
$ws.Range("J137").Value = 5600
$ws.Range("K137").Value = 100126710
$ws.Range("L137").Value = 16800
$ws.Range("M137").Value = -100124160
$ws.Range("N137").Value = -21900

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 50000
$ws.Range("J44").Value = 50000
$ws.Range("L44").Value = 50000
$ws.Range("N44").Value = -50976
$ws.Range("H63").Value = 3710.8572
$ws.Range("I63").Value = 3116.16
$ws.Range("J63").Value = 8666.666999999999
$ws.Range("K63").Value = 3116.16
$ws.Range("L63").Value = 8666.666999999999
$ws.Range("M63").Value = -2430.16
$ws.Range("N63").Value = -10038.667
$ws.Range("H66").Value = 3710.8572
$ws.Range("I66").Value = 3116.16
$ws.Range("J66").Value = 8666.666999999999
$ws.Range("K66").Value = 15580.8
$ws.Range("L66").Value = 43333.335
$ws.Range("M66").Value = -12148.8
$ws.Range("N66").Value = -50197.335
$ws.Range("H80").Value = 48419.5
$ws.Range("J80").Value = 48419.5
$ws.Range("L80").Value = 48419.5
$ws.Range("N80").Value = -50415.5
$ws.Range("H83").Value = 48419.5
$ws.Range("J83").Value = 48419.5
$ws.Range("L83").Value = 145258.5
$ws.Range("N83").Value = -155242.5
$ws.Range("H110").Value = 2286.8147
$ws.Range("I110").Value = 2578.2104
$ws.Range("K110").Value = 2578.2104
$ws.Range("M110").Value = -533.2103999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 34902.418
$ws.Range("I20").Value = 1874.1578
$ws.Range("J20").Value = 87197.164
$ws.Range("K20").Value = 1874.1578
$ws.Range("L20").Value = 87197.164
$ws.Range("M20").Value = -1627.1578
$ws.Range("N20").Value = -87691.164
$ws.Range("H54").Value = 1992.4286
$ws.Range("I54").Value = 2074.5
$ws.Range("K54").Value = 2074.5
$ws.Range("M54").Value = -1590.5
$ws.Range("H134").Value = 4477.8335
$ws.Range("I134").Value = 4477.8335
$ws.Range("K134").Value = 13433.5005
$ws.Range("M134").Value = -10898.5005

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3552.9443
$ws.Range("I16").Value = 1288.9166
$ws.Range("J16").Value = 8081
$ws.Range("K16").Value = 1288.9166
$ws.Range("L16").Value = 8081
$ws.Range("M16").Value = -1001.9166
$ws.Range("N16").Value = -8655
$ws.Range("H22").Value = 2275.3635
$ws.Range("I22").Value = 3257.1428
$ws.Range("K22").Value = 3257.1428
$ws.Range("M22").Value = -2907.1428
$ws.Range("H31").Value = 4072.4707
$ws.Range("I31").Value = 4363.6284
$ws.Range("J31").Value = 3435.5625
$ws.Range("K31").Value = 4363.6284
$ws.Range("L31").Value = 3435.5625
$ws.Range("M31").Value = -4068.6284
$ws.Range("N31").Value = -4025.5625
$ws.Range("H34").Value = 4072.4707
$ws.Range("I34").Value = 4363.6284
$ws.Range("J34").Value = 3435.5625
$ws.Range("K34").Value = 4363.6284
$ws.Range("L34").Value = 3435.5625
$ws.Range("M34").Value = -4161.6284
$ws.Range("N34").Value = -3839.5625
$ws.Range("H58").Value = 3584.8235
$ws.Range("J58").Value = 4108
$ws.Range("L58").Value = 4108
$ws.Range("N58").Value = -4514
$ws.Range("H113").Value = 3552.9443
$ws.Range("I113").Value = 1288.9166
$ws.Range("J113").Value = 8081
$ws.Range("K113").Value = 1288.9166
$ws.Range("L113").Value = 8081
$ws.Range("M113").Value = 881.0834
$ws.Range("N113").Value = -12421
$ws.Range("H136").Value = 3584.8235
$ws.Range("J136").Value = 4108
$ws.Range("L136").Value = 12324
$ws.Range("N136").Value = -17424

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 937.4545000000001
$ws.Range("I8").Value = 937.4545000000001
$ws.Range("K8").Value = 2812.3635
$ws.Range("M8").Value = -2673.3635
$ws.Range("H57").Value = 200124.5
$ws.Range("I57").Value = 800
$ws.Range("J57").Value = 399449
$ws.Range("K57").Value = 2400
$ws.Range("L57").Value = 1198347
$ws.Range("M57").Value = -1841
$ws.Range("N57").Value = -1199465
$ws.Range("H68").Value = 1690.6666
$ws.Range("J68").Value = 2222
$ws.Range("L68").Value = 6666
$ws.Range("N68").Value = -8288
$ws.Range("H71").Value = 1690.6666
$ws.Range("J71").Value = 2222
$ws.Range("L71").Value = 19998
$ws.Range("N71").Value = -28110
$ws.Range("H112").Value = 1200
$ws.Range("I112").Value = 1200
$ws.Range("K112").Value = 3600
$ws.Range("M112").Value = -2492
$ws.Range("H129").Value = 1264.8334
$ws.Range("J129").Value = 2113
$ws.Range("L129").Value = 6339
$ws.Range("N129").Value = -16339

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 141763.75
$ws.Range("I80").Value = 161301.58
$ws.Range("J80").Value = 4999
$ws.Range("K80").Value = 161301.58
$ws.Range("L80").Value = 4999
$ws.Range("M80").Value = -160303.58
$ws.Range("N80").Value = -6995
$ws.Range("H83").Value = 141763.75
$ws.Range("I83").Value = 161301.58
$ws.Range("J83").Value = 4999
$ws.Range("K83").Value = 806507.8999999999
$ws.Range("L83").Value = 24995
$ws.Range("M83").Value = -801515.8999999999
$ws.Range("N83").Value = -34979
$ws.Range("H97").Value = 34833.332
$ws.Range("I97").Value = 1500
$ws.Range("J97").Value = 51500
$ws.Range("K97").Value = 1500
$ws.Range("L97").Value = 51500
$ws.Range("M97").Value = -1004
$ws.Range("N97").Value = -52492
$ws.Range("H99").Value = 15999.75
$ws.Range("I99").Value = 8333
$ws.Range("J99").Value = 39000
$ws.Range("K99").Value = 8333
$ws.Range("L99").Value = 39000
$ws.Range("M99").Value = -6087
$ws.Range("N99").Value = -43492
$ws.Range("H113").Value = 11518.066
$ws.Range("I113").Value = 6985
$ws.Range("J113").Value = 15484.5
$ws.Range("K113").Value = 6985
$ws.Range("L113").Value = 15484.5
$ws.Range("M113").Value = -4815
$ws.Range("N113").Value = -19824.5
$ws.Range("H132").Value = 9306.075000000001
$ws.Range("I132").Value = 8728.615
$ws.Range("J132").Value = 10378.5
$ws.Range("K132").Value = 26185.845
$ws.Range("L132").Value = 31135.5
$ws.Range("M132").Value = -23655.845
$ws.Range("N132").Value = -36195.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 99.8
$ws.Range("I22").Value = 99
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 99
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = 196
$ws.Range("N22").Value = -690
$ws.Range("H27").Value = 99.8
$ws.Range("I27").Value = 99
$ws.Range("J27").Value = 100
$ws.Range("K27").Value = 99
$ws.Range("L27").Value = 100
$ws.Range("M27").Value = 8
$ws.Range("N27").Value = -314
$ws.Range("H46").Value = 3701.3547
$ws.Range("J46").Value = 3808.0667
$ws.Range("L46").Value = 3808.0667
$ws.Range("N46").Value = -4184.066699999999
$ws.Range("H55").Value = 775.1818
$ws.Range("I55").Value = 552.7
$ws.Range("K55").Value = 552.7
$ws.Range("M55").Value = -379.7
$ws.Range("H132").Value = 4688.0347
$ws.Range("I132").Value = 4782.32
$ws.Range("K132").Value = 14346.96
$ws.Range("M132").Value = -11816.96

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 22713
$ws.Range("J41").Value = 22713
$ws.Range("L41").Value = 22713
$ws.Range("N41").Value = -23493
$ws.Range("H45").Value = 27050
$ws.Range("J45").Value = 27050
$ws.Range("L45").Value = 27050
$ws.Range("N45").Value = -28032
$ws.Range("H46").Value = 61710
$ws.Range("J46").Value = 61710
$ws.Range("L46").Value = 61710
$ws.Range("N46").Value = -62172
$ws.Range("H86").Value = 66706.25
$ws.Range("J86").Value = 66706.25
$ws.Range("L86").Value = 66706.25
$ws.Range("N86").Value = -68952.25
$ws.Range("H89").Value = 66706.25
$ws.Range("J89").Value = 66706.25
$ws.Range("L89").Value = 333531.25
$ws.Range("N89").Value = -344763.25
$ws.Range("H134").Value = 61710
$ws.Range("J134").Value = 61710
$ws.Range("L134").Value = 185130
$ws.Range("N134").Value = -190200

Write-Output "Applied 249 cell edits across 8 sheets"
